# Tweaker to reduce execution speed for _X strategies:
# - Remove the two extra ScalpEmaRsiAdx_X test rows (rows 3 and 4), shifting
#   the rest of the (empty) template rows up.
# - Update the remaining ScalpEmaRsiAdx_X test row's date range and interval.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete rows 3 and 4 entirely (shifts everything below up by two rows).
$ws.Range("A3:A4").EntireRow.Delete()

# Update row 2: widen the "From" date and switch interval to 3m.
$ws.Range("D2").Value = 44197
$ws.Range("F2").Value = "3m"

# Update the active selection to D3.
[void]$ws.Range("D3").Select()
